# Edit script: rename stat sheets to human-friendly titles and bump the
# "age in days" component of every player's Age column (E) by one day
# across all per-player stats sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets (tab names only; sheetId / r:id stay the same) ---
$renameMap = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"    = "Shooting Stats"
    "PassingStats"     = "Passing Stats"
    "PassTypes"        = "Pass Types"
    "GoalShotCreation" = "Goal & Shot Creation"
    "DefensiveActions" = "Defensive Actions"
    "PlayingTime"      = "Playing Time"
    "MiscStats"        = "Miscellaneous Stats"
}

foreach ($oldName in $renameMap.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renameMap[$oldName]
}

# --- 2. Bump the day component of every "YY-DDD" age value in column E ---
# Applies to every stats sheet (everything except "Matches"); player rows
# run from row 4 down to the last row whose E value still matches the
# "YY-DDD" pattern (summary rows below use a decimal, e.g. "27.3", and are
# left untouched).
$ageRegex = [regex]"^(\d+)-(\d{3})$"

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Matches") {
        continue
    }

    $row = 4
    while ($true) {
        $cell = $ws.Cells.Item($row, 5)
        $val = $cell.Value2
        if ($val -eq $null) {
            break
        }

        $text = [string]$val
        $m = $ageRegex.Match($text)
        if (-not $m.Success) {
            break
        }

        $years = $m.Groups[1].Value
        $days = [int]$m.Groups[2].Value + 1
        $newVal = "{0}-{1:D3}" -f $years, $days
        $cell.Value = $newVal

        $row = $row + 1
    }
}
